$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-01-04 Saturday", $true, $true, $false, $false, $false, $true, 1, $false, "2025-01-05 Sunday", 2) | Out-Null
$d.Content.Find.Execute("22+67=89", $true, $true, $false, $false, $false, $true, 1, $false, "96-10=86", 2) | Out-Null
$d.Content.Find.Execute("10+75=85", $true, $true, $false, $false, $false, $true, 1, $false, "2+61=63", 2) | Out-Null
$d.Content.Find.Execute("52-37=15", $true, $true, $false, $false, $false, $true, 1, $false, "7+57=64", 2) | Out-Null
$d.Content.Find.Execute("72-65=7", $true, $true, $false, $false, $false, $true, 1, $false, "15+43=58", 2) | Out-Null
$d.Content.Find.Execute("8+2=10", $true, $true, $false, $false, $false, $true, 1, $false, "3+71=74", 2) | Out-Null
$d.Content.Find.Execute("41-29=12", $true, $true, $false, $false, $false, $true, 1, $false, "51-8=43", 2) | Out-Null
$d.Content.Find.Execute("12+3=15", $true, $true, $false, $false, $false, $true, 1, $false, "48-29=19", 2) | Out-Null
$d.Content.Find.Execute("13+27=40", $true, $true, $false, $false, $false, $true, 1, $false, "6+8=14", 2) | Out-Null
$d.Content.Find.Execute("86-58=28", $true, $true, $false, $false, $false, $true, 1, $false, "57+15=72", 2) | Out-Null
$d.Content.Find.Execute("62-16=46", $true, $true, $false, $false, $false, $true, 1, $false, "58-33=25", 2) | Out-Null
$d.Content.Find.Execute("37-6=31", $true, $true, $false, $false, $false, $true, 1, $false, "77+21=98", 2) | Out-Null
$d.Content.Find.Execute("80-57=23", $true, $true, $false, $false, $false, $true, 1, $false, "60-17=43", 2) | Out-Null
$d.Content.Find.Execute("88-70=18", $true, $true, $false, $false, $false, $true, 1, $false, "21+35=56", 2) | Out-Null
$d.Content.Find.Execute("72-17=55", $true, $true, $false, $false, $false, $true, 1, $false, "31+47=78", 2) | Out-Null
$d.Content.Find.Execute("5+28=33", $true, $true, $false, $false, $false, $true, 1, $false, "29-16=13", 2) | Out-Null
$d.Content.Find.Execute("28+56=84", $true, $true, $false, $false, $false, $true, 1, $false, "91-61=30", 2) | Out-Null
$d.Content.Find.Execute("80-73=7", $true, $true, $false, $false, $false, $true, 1, $false, "29+64=93", 2) | Out-Null
$d.Content.Find.Execute("73-26=47", $true, $true, $false, $false, $false, $true, 1, $false, "77-33=44", 2) | Out-Null
$d.Content.Find.Execute("69+3=72", $true, $true, $false, $false, $false, $true, 1, $false, "12+57=69", 2) | Out-Null
$d.Content.Find.Execute("85-12=73", $true, $true, $false, $false, $false, $true, 1, $false, "81-3=78", 2) | Out-Null
$d.Content.Find.Execute("85-14=71", $true, $true, $false, $false, $false, $true, 1, $false, "46+52=98", 2) | Out-Null
$d.Content.Find.Execute("28+28=56", $true, $true, $false, $false, $false, $true, 1, $false, "13+75=88", 2) | Out-Null
$d.Content.Find.Execute("99-77=22", $true, $true, $false, $false, $false, $true, 1, $false, "28-2=26", 2) | Out-Null
$d.Content.Find.Execute("15+17=32", $true, $true, $false, $false, $false, $true, 1, $false, "3+16=19", 2) | Out-Null
$d.Content.Find.Execute("71+2=73", $true, $true, $false, $false, $false, $true, 1, $false, "56+35=91", 2) | Out-Null
$d.Content.Find.Execute("75-29=46", $true, $true, $false, $false, $false, $true, 1, $false, "91-0=91", 2) | Out-Null
$d.Content.Find.Execute("9+31=40", $true, $true, $false, $false, $false, $true, 1, $false, "34+16=50", 2) | Out-Null
$d.Content.Find.Execute("28+36=64", $true, $true, $false, $false, $false, $true, 1, $false, "36-20=16", 2) | Out-Null
$d.Content.Find.Execute("92-68=24", $true, $true, $false, $false, $false, $true, 1, $false, "2+36=38", 2) | Out-Null
$d.Content.Find.Execute("85-6=79", $true, $true, $false, $false, $false, $true, 1, $false, "52-29=23", 2) | Out-Null
$d.Content.Find.Execute("25+70=95", $true, $true, $false, $false, $false, $true, 1, $false, "86-60=26", 2) | Out-Null
$d.Content.Find.Execute("0+86=86", $true, $true, $false, $false, $false, $true, 1, $false, "44-4=40", 2) | Out-Null
$d.Content.Find.Execute("20+79=99", $true, $true, $false, $false, $false, $true, 1, $false, "84+0=84", 2) | Out-Null
$d.Content.Find.Execute("22-13=9", $true, $true, $false, $false, $false, $true, 1, $false, "55+35=90", 2) | Out-Null
$d.Content.Find.Execute("0+71=71", $true, $true, $false, $false, $false, $true, 1, $false, "90-60=30", 2) | Out-Null
$d.Content.Find.Execute("79-19=60", $true, $true, $false, $false, $false, $true, 1, $false, "57-31=26", 2) | Out-Null
$d.Content.Find.Execute("5+29=34", $true, $true, $false, $false, $false, $true, 1, $false, "31-25=6", 2) | Out-Null
$d.Content.Find.Execute("24-11=13", $true, $true, $false, $false, $false, $true, 1, $false, "74-53=21", 2) | Out-Null
$d.Content.Find.Execute("45-11=34", $true, $true, $false, $false, $false, $true, 1, $false, "88-80=8", 2) | Out-Null
$d.Content.Find.Execute("20+17=37", $true, $true, $false, $false, $false, $true, 1, $false, "37+28=65", 2) | Out-Null
$d.Content.Find.Execute("6+53=59", $true, $true, $false, $false, $false, $true, 1, $false, "50+3=53", 2) | Out-Null
$d.Content.Find.Execute("91-82=9", $true, $true, $false, $false, $false, $true, 1, $false, "0+4=4", 2) | Out-Null
$d.Content.Find.Execute("73-7=66", $true, $true, $false, $false, $false, $true, 1, $false, "73-70=3", 2) | Out-Null
$d.Content.Find.Execute("9+71=80", $true, $true, $false, $false, $false, $true, 1, $false, "6+24=30", 2) | Out-Null
$d.Content.Find.Execute("24+75=99", $true, $true, $false, $false, $false, $true, 1, $false, "60-56=4", 2) | Out-Null
$d.Content.Find.Execute("92-60=32", $true, $true, $false, $false, $false, $true, 1, $false, "26+24=50", 2) | Out-Null
$d.Content.Find.Execute("75-65=10", $true, $true, $false, $false, $false, $true, 1, $false, "87-23=64", 2) | Out-Null
$d.Content.Find.Execute("50+25=75", $true, $true, $false, $false, $false, $true, 1, $false, "49-37=12", 2) | Out-Null
$d.Content.Find.Execute("98-80=18", $true, $true, $false, $false, $false, $true, 1, $false, "18+27=45", 2) | Out-Null
$d.Content.Find.Execute("9+64=73", $true, $true, $false, $false, $false, $true, 1, $false, "97-57=40", 2) | Out-Null
$d.Content.Find.Execute("4+63=67", $true, $true, $false, $false, $false, $true, 1, $false, "39+18=57", 2) | Out-Null
$d.Content.Find.Execute("27+38=65", $true, $true, $false, $false, $false, $true, 1, $false, "60+36=96", 2) | Out-Null
$d.Content.Find.Execute("24+1=25", $true, $true, $false, $false, $false, $true, 1, $false, "40+32=72", 2) | Out-Null
$d.Content.Find.Execute("54-16=38", $true, $true, $false, $false, $false, $true, 1, $false, "7+75=82", 2) | Out-Null
$d.Content.Find.Execute("86-27=59", $true, $true, $false, $false, $false, $true, 1, $false, "52+28=80", 2) | Out-Null
$d.Content.Find.Execute("28-5=23", $true, $true, $false, $false, $false, $true, 1, $false, "98-51=47", 2) | Out-Null
$d.Content.Find.Execute("53+21=74", $true, $true, $false, $false, $false, $true, 1, $false, "74+13=87", 2) | Out-Null
$d.Content.Find.Execute("60-11=49", $true, $true, $false, $false, $false, $true, 1, $false, "31+13=44", 2) | Out-Null
$d.Content.Find.Execute("11-5=6", $true, $true, $false, $false, $false, $true, 1, $false, "83-11=72", 2) | Out-Null
$d.Content.Find.Execute("58-15=43", $true, $true, $false, $false, $false, $true, 1, $false, "80-66=14", 2) | Out-Null
$d.Content.Find.Execute("31-27=4", $true, $true, $false, $false, $false, $true, 1, $false, "60+21=81", 2) | Out-Null
$d.Content.Find.Execute("17+42=59", $true, $true, $false, $false, $false, $true, 1, $false, "9+50=59", 2) | Out-Null
$d.Content.Find.Execute("72+25=97", $true, $true, $false, $false, $false, $true, 1, $false, "26+51=77", 2) | Out-Null
$d.Content.Find.Execute("26-19=7", $true, $true, $false, $false, $false, $true, 1, $false, "44-2=42", 2) | Out-Null
$d.Content.Find.Execute("12+37=49", $true, $true, $false, $false, $false, $true, 1, $false, "59+2=61", 2) | Out-Null
$d.Content.Find.Execute("78-60=18", $true, $true, $false, $false, $false, $true, 1, $false, "98-47=51", 2) | Out-Null
$d.Content.Find.Execute("2+76=78", $true, $true, $false, $false, $false, $true, 1, $false, "84-70=14", 2) | Out-Null
$d.Content.Find.Execute("62-21=41", $true, $true, $false, $false, $false, $true, 1, $false, "78-77=1", 2) | Out-Null
$d.Content.Find.Execute("73-55=18", $true, $true, $false, $false, $false, $true, 1, $false, "63-1=62", 2) | Out-Null
$d.Content.Find.Execute("27-3=24", $true, $true, $false, $false, $false, $true, 1, $false, "81-55=26", 2) | Out-Null
$d.Content.Find.Execute("42-36=6", $true, $true, $false, $false, $false, $true, 1, $false, "82-44=38", 2) | Out-Null
$d.Content.Find.Execute("59-4=55", $true, $true, $false, $false, $false, $true, 1, $false, "89-17=72", 2) | Out-Null
$d.Content.Find.Execute("94-37=57", $true, $true, $false, $false, $false, $true, 1, $false, "73-38=35", 2) | Out-Null
$d.Content.Find.Execute("97-50=47", $true, $true, $false, $false, $false, $true, 1, $false, "83-42=41", 2) | Out-Null
$d.Content.Find.Execute("81-25=56", $true, $true, $false, $false, $false, $true, 1, $false, "38-9=29", 2) | Out-Null
$d.Content.Find.Execute("83-65=18", $true, $true, $false, $false, $false, $true, 1, $false, "98-35=63", 2) | Out-Null
$d.Content.Find.Execute("18+32=50", $true, $true, $false, $false, $false, $true, 1, $false, "96-14=82", 2) | Out-Null
$d.Content.Find.Execute("44-6=38", $true, $true, $false, $false, $false, $true, 1, $false, "95-7=88", 2) | Out-Null
$d.Content.Find.Execute("21+13=34", $true, $true, $false, $false, $false, $true, 1, $false, "69+19=88", 2) | Out-Null
$d.Content.Find.Execute("29+11=40", $true, $true, $false, $false, $false, $true, 1, $false, "77-23=54", 2) | Out-Null
$d.Content.Find.Execute("11+71=82", $true, $true, $false, $false, $false, $true, 1, $false, "68+8=76", 2) | Out-Null
$d.Content.Find.Execute("81+1=82", $true, $true, $false, $false, $false, $true, 1, $false, "24+59=83", 2) | Out-Null
$d.Content.Find.Execute("90-67=23", $true, $true, $false, $false, $false, $true, 1, $false, "10+72=82", 2) | Out-Null
$d.Content.Find.Execute("39-17=22", $true, $true, $false, $false, $false, $true, 1, $false, "95-73=22", 2) | Out-Null
$d.Content.Find.Execute("71+12=83", $true, $true, $false, $false, $false, $true, 1, $false, "7+27=34", 2) | Out-Null
$d.Content.Find.Execute("86-21=65", $true, $true, $false, $false, $false, $true, 1, $false, "45+0=45", 2) | Out-Null
$d.Content.Find.Execute("66-10=56", $true, $true, $false, $false, $false, $true, 1, $false, "99-45=54", 2) | Out-Null
$d.Content.Find.Execute("40+49=89", $true, $true, $false, $false, $false, $true, 1, $false, "18+38=56", 2) | Out-Null
$d.Content.Find.Execute("38+5=43", $true, $true, $false, $false, $false, $true, 1, $false, "22+75=97", 2) | Out-Null
$d.Content.Find.Execute("25+1=26", $true, $true, $false, $false, $false, $true, 1, $false, "10+32=42", 2) | Out-Null
$d.Content.Find.Execute("0+82=82", $true, $true, $false, $false, $false, $true, 1, $false, "78+13=91", 2) | Out-Null
$d.Content.Find.Execute("45-42=3", $true, $true, $false, $false, $false, $true, 1, $false, "18+61=79", 2) | Out-Null
$d.Content.Find.Execute("50-42=8", $true, $true, $false, $false, $false, $true, 1, $false, "17-10=7", 2) | Out-Null
$d.Content.Find.Execute("49+31=80", $true, $true, $false, $false, $false, $true, 1, $false, "94-72=22", 2) | Out-Null
$d.Content.Find.Execute("87-51=36", $true, $true, $false, $false, $false, $true, 1, $false, "97-27=70", 2) | Out-Null
$d.Content.Find.Execute("72-37=35", $true, $true, $false, $false, $false, $true, 1, $false, "98-67=31", 2) | Out-Null
$d.Content.Find.Execute("78-19=59", $true, $true, $false, $false, $false, $true, 1, $false, "95-73=22", 2) | Out-Null
$d.Content.Find.Execute("18+40=58", $true, $true, $false, $false, $false, $true, 1, $false, "90-62=28", 2) | Out-Null
$d.Content.Find.Execute("79-7=72", $true, $true, $false, $false, $false, $true, 1, $false, "47+5=52", 2) | Out-Null
$d.Content.Find.Execute("62-39=23", $true, $true, $false, $false, $false, $true, 1, $false, "69+14=83", 2) | Out-Null
